$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Experimental flag: "false" -> "true" (row 7) ---
# --- Case Sensitive value: blank -> "true" (row 14) ---
# Plain .Value assignment of the literal text "true"/"false" gets auto-coerced
# to a boolean by Excel, so first mirror the existing text-typed "false"
# cell into B14 (Copy preserves the Text type), then overwrite both cells'
# text via a scratch formula cell + PasteSpecial(values), which also keeps
# the Text type instead of re-parsing into a boolean.
$ws.Range("B7").Copy($ws.Range("B14")) | Out-Null

$scratch = $ws.Range("Z1")
$scratch.Formula = "=""true"""
$scratch.Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4163) | Out-Null
$ws.Range("B14").PasteSpecial(-4163) | Out-Null
$scratch.Value = ""

# --- Date updated ---
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
